# Weekly update: insert a new price record as row 32 (Comercializadora del
# Agro de Limarí - Zapallo italiano), pushing the previously existing rows
# 32-45 down to 33-46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 32, shifting rows 32:45 down to 33:46.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record.
$ws.Cells.Item(32, 1).Value = 2
$ws.Cells.Item(32, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(32, 3).Value = "Coquimbo"
$ws.Cells.Item(32, 4).Value = 44567
$ws.Cells.Item(32, 5).Value = 4
$ws.Cells.Item(32, 6).Value = 100112032
$ws.Cells.Item(32, 7).Value = "Zapallo italiano"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 400
$ws.Cells.Item(32, 11).Value = 6500
$ws.Cells.Item(32, 12).Value = 7000
$ws.Cells.Item(32, 13).Value = 6750
$ws.Cells.Item(32, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(32, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(32, 16).Value = 112
$ws.Cells.Item(32, 17).Value = 60
$ws.Cells.Item(32, 18).Value = "Hortaliza"

# Match the date cell style used by the rest of column D (numeric date format).
$ws.Cells.Item(32, 4).NumberFormat = $ws.Cells.Item(33, 4).NumberFormat
